$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.468.86"
$ws.Range("E2").Value = "  -2.65%  "
$ws.Range("D3").Value = "1.805.33"
$ws.Range("E3").Value = "  -2.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.72%  "
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "308.35"
$ws.Range("E6").Value = "  -1.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4540"
$ws.Range("E7").Value = "  -1.42%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3655"
$ws.Range("E8").Value = "  -1.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07115"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8749"
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07735"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.33"
$ws.Range("E12").Value = "  -3.37%  "
$ws.Range("D13").Value = "1.830.56"
$ws.Range("E13").Value = "  -5.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.268"
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.340"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "85.92"
$ws.Range("E16").Value = "  -5.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008572"
$ws.Range("E18").Value = "  -3.93%  "
$ws.Range("D20").Value = "26.517.58"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.22"
$ws.Range("E21").Value = "  -3.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.969"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.974"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.50"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("E26").Value = "  -2.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.997"
$ws.Range("E27").Value = "  -3.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "112.47"
$ws.Range("E28").Value = "  -2.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.813"
$ws.Range("E29").Value = "  -4.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08652"
$ws.Range("E30").Value = "  -1.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.048"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7278"
$ws.Range("E32").Value = "  -4.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.425"
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("E34").Value = "  -4.75%  "
$ws.Range("E35").Value = "  +0.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.554"
$ws.Range("E36").Value = "  -6.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.080"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01927"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05086"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.880"
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.933"
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4991"
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1562"
$ws.Range("E43").Value = "  -3.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.095"
$ws.Range("E44").Value = "  -3.30%  "
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("E46").Value = "  -3.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.946"
$ws.Range("E47").Value = "  -3.03%  "
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.585"
$ws.Range("E49").Value = "  -2.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05985"
$ws.Range("E50").Value = "  -3.65%  "
$ws.Range("E51").Value = "  -2.66%  "
